$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.840.17'
$ws.Range("E2").Value = '  +0.87%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.809.29'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.92'
$ws.Range("E5").Value = '  +2.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.607'
$ws.Range("E6").Value = '  +0.59%  '
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.50'
$ws.Range("E8").Value = '  -5.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.309'
$ws.Range("E9").Value = '  +5.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0680'
$ws.Range("E10").Value = '  +2.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0998'
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.071.05'
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.809.11'
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.01'
$ws.Range("E14").Value = '  +0.82%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.654'
$ws.Range("E15").Value = '  +4.13%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.63'
$ws.Range("E16").Value = '  +5.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '34.825.50'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.63'
$ws.Range("E18").Value = '  +1.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0780'
$ws.Range("E19").Value = '  +1.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '236.47'
$ws.Range("E20").Value = '  -1.72%  '
$ws.Range("E21").Value = '  +4.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.64'
$ws.Range("E22").Value = '  +7.14%  '
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.26'
$ws.Range("E24").Value = '  +4.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '173.14'
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.74'
$ws.Range("E26").Value = '  +1.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.30'
$ws.Range("E27").Value = '  -0.73%  '
$ws.Range("E28").Value = '  -1.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.58'
$ws.Range("E29").Value = '  +29.12%  '
$ws.Range("E30").Value = '  +0.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.339.46'
$ws.Range("E31").Value = '  +37.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0546'
$ws.Range("E32").Value = '  +6.48%  '
$ws.Range("E33").Value = '  +1.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.94'
$ws.Range("E34").Value = '  +1.68%  '
$ws.Range("E35").Value = '  -0.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.17'
$ws.Range("E36").Value = '  +10.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '93.15'
$ws.Range("E37").Value = '  +6.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.673'
$ws.Range("E38").Value = '  +3.89%  '
$ws.Range("E39").Value = '  +1.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.301.87'
$ws.Range("E40").Value = '  -1.11%  '
$ws.Range("E41").Value = '  +4.04%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.980'
$ws.Range("E42").Value = '  +4.42%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.62'
$ws.Range("E43").Value = '  -0.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.31'
$ws.Range("E44").Value = '  -1.11%  '
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("E46").Value = '  -1.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.23'
$ws.Range("E47").Value = '  +7.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0511'
$ws.Range("E48").Value = '  -1.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.985.19'
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("E50").Value = '  +0.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0642'
$ws.Range("E51").Value = '  +5.76%  '
